$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colors")

# Insert a new row at 15 (shifts old rows 15-27 down to 16-28)
$ws.Rows.Item(15).Insert()

# Set new row 15 cell values
$ws.Range("A15").Value = "Thumbnail disc"
$ws.Range("B15").Value = "Text"
$ws.Range("B15").Font.Color = 8355711
$ws.Range("B15").Interior.Color = 13231601
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("F15").Value = "Thumbnail disc for initials"

# New cell in row 2
$ws.Range("L2").Value = "Background"

# Update F9 text (create shared string "Popup menus background" before "Popup menus title")
$ws.Range("F9").Value = "Popup menus background"

# New cell in row 5
$ws.Range("F5").Value = "Popup menus title"
